# manufacture/PurchaseList.xlsx -- "update pcb and submit pcb file to betalayout"
#
# Refresh the Supplier Stock (J), Adjusted Supplier Order Qty (N) and Supplier
# Unit Price (P) figures pulled from the distributor for every BOM line, plus the
# "report created" time stamp (F24) and the production-run size (N25) shown in
# the footer. O (ROUNDUP) and Q (price * qty) are formulas and recompute on their
# own, as do the SUM/NOW()-driven totals in row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer labels are entered with a leading apostrophe in the sheet so Excel keeps
# them as text (quote-prefixed) even though they look numeric.
$ws.Range("F24").Value = "'16:31"   # report created: <time>
$ws.Range("N25").Value = "'36"      # production QTY:

# Row 2 (C1, C4, C5)
$ws.Range("J2").Value = 37050
$ws.Range("N2").Value = 108
$ws.Range("P2").Value = 0.019

# Row 3 (C2, C3)
$ws.Range("N3").Value = 100
$ws.Range("P3").Value = 0.0362

# Row 4 (C6, C7)
$ws.Range("J4").Value = 33872
$ws.Range("N4").Value = 72
$ws.Range("P4").Value = 0.01448

# Row 5 (C8, C10)
$ws.Range("J5").Value = 9999
$ws.Range("N5").Value = 100
$ws.Range("P5").Value = 0.02896

# Row 6 (D1)
$ws.Range("J6").Value = 59043
$ws.Range("P6").Value = 0.06877

# Row 7 (D2)
$ws.Range("J7").Value = 98272
$ws.Range("N7").Value = 36
$ws.Range("P7").Value = 0.10497

# Row 8 (D3)
$ws.Range("J8").Value = 57880
$ws.Range("N8").Value = 36
$ws.Range("P8").Value = 0.27147

# Row 9 (D5)
$ws.Range("J9").Value = 116092
$ws.Range("N9").Value = 36
$ws.Range("P9").Value = 0.34477

# Row 10 (JP2)
$ws.Range("J10").Value = 21003
$ws.Range("N10").Value = 36
$ws.Range("P10").Value = 0.67325

# Row 11 (Q1)
$ws.Range("J11").Value = 94764
$ws.Range("N11").Value = 36
$ws.Range("P11").Value = 0.2407

# Row 12 (R1, R3, R4, R9, R14)
$ws.Range("J12").Value = 1836811
$ws.Range("N12").Value = 180
$ws.Range("P12").Value = 0.01176

# Row 13 (R2)
$ws.Range("J13").Value = 68380
$ws.Range("N13").Value = 36
$ws.Range("P13").Value = 0.02896

# Row 14 (R5, R7)
$ws.Range("J14").Value = 83930
$ws.Range("N14").Value = 100
$ws.Range("P14").Value = 0.01176

# Row 15 (R6, R8)
$ws.Range("J15").Value = 131942
$ws.Range("N15").Value = 100
$ws.Range("P15").Value = 0.01176

# Row 16 (R10)
$ws.Range("J16").Value = 98927
$ws.Range("N16").Value = 36
$ws.Range("P16").Value = 0.02896

# Row 17 (R12)
$ws.Range("J17").Value = 891522
$ws.Range("N17").Value = 36
$ws.Range("P17").Value = 0.01448

# Row 18 (R15)
$ws.Range("J18").Value = 16412
$ws.Range("N18").Value = 36
$ws.Range("P18").Value = 0.03529

# Row 19 (SW3, SW4)
$ws.Range("J19").Value = 14994
$ws.Range("N19").Value = 72
$ws.Range("P19").Value = 0.65244

# Row 20 (U1)
$ws.Range("N20").Value = 36
$ws.Range("P20").Value = 3.47

# Row 21 (U2)
$ws.Range("J21").Value = 27364
$ws.Range("N21").Value = 36
$ws.Range("P21").Value = 0.60629

# Row 22 (Y1)
$ws.Range("J22").Value = 530
$ws.Range("N22").Value = 36
$ws.Range("P22").Value = 1.23
